# 成品糖.xlsx edit:
#  1. For every year block (4 rows: A/B/C/D sub-periods), swap the 2nd and
#     3rd rows (the "B" and "C" sub-period rows) — their A:E contents trade
#     places.
#  2. Delete columns F ("成品糖产销率") and G ("成品糖销售量") entirely, since
#     that data is dropped from the sheet (dimension shrinks from A1:G65 to
#     A1:E65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data starts at row 2 and is laid out in 4-row blocks (one per year):
# rowA, rowB, rowC, rowD. Swap rowB <-> rowC within each block.
$firstDataRow = 2
$blockSize = 4
$lastDataRow = 65

for ($rowB = $firstDataRow + 1; $rowB -le $lastDataRow; $rowB += $blockSize) {
    $rowC = $rowB + 1

    $valsB = $ws.Range("A$rowB`:E$rowB").Value2
    $valsC = $ws.Range("A$rowC`:E$rowC").Value2

    $ws.Range("A$rowB`:E$rowB").Value2 = $valsC
    $ws.Range("A$rowC`:E$rowC").Value2 = $valsB
}

# Drop columns F and G completely.
$ws.Range("F1:G1").EntireColumn.Delete()
